$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 832 (the "「ゆっくりと」بشويش" entry), shifting rows below it up.
$ws.Rows.Item(832).Delete()
